$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the value currently in A414 ("Wikipedia:Protection policy") up to A393,
# shifting the existing values in A393:A413 down by one row into A394:A414.

$startRow = 393
$endRow = 414

# Capture the value that needs to move to the top of the block.
$movingValue = $ws.Cells.Item($endRow, 1).Value2

# Shift rows down, working from the bottom of the block upward so we don't
# overwrite values before they've been read.
for ($r = $endRow; $r -gt $startRow; $r--) {
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r - 1, 1).Value2
}

# Place the moved value at the top of the block.
$ws.Cells.Item($startRow, 1).Value2 = $movingValue
